# "Atualiza copa leon e reorganiza copa"
# Insert a new team ("GrioTeam", id 14933455) right after "GE Bebum" (row 9),
# pushing "Grêmio_Campeão_LA_27" and every team below it down by one row,
# and rebuild the hyperlinks in column C so they keep pointing at the right rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10 - this shifts rows 10:20 down to 11:21
# (values, number formats and styles all move down automatically).
$ws.Rows("10:10").Insert()

# Fill in the new row for "GrioTeam".
$ws.Range("A10").Value2 = "GrioTeam"
$ws.Range("B10").Value2 = 14933455
$ws.Range("C10").Value2 = "https://cartola.globo.com/#!/time/14933455"

# Row-insert does not renumber the worksheet's <hyperlinks> collection, so every
# hyperlink from the old row 10 downward now points one row too high. Rebuild
# the whole C2:C21 hyperlink list from scratch, in final row order.
$ws.Hyperlinks.Delete()

$teamIds = @(
    19833277,
    19209079,
    1488983,
    287965,
    2916559,
    186283,
    2371918,
    16411206,
    14933455,
    47775950,
    1747619,
    32966,
    44810918,
    1867254,
    4088673,
    1326835,
    20651178,
    14709358,
    184499,
    1273719
)

for ($i = 0; $i -lt $teamIds.Length; $i++) {
    $row = 2 + $i
    $id = $teamIds[$i]
    $ws.Hyperlinks.Add($ws.Range("C$row"), "https://cartola.globo.com/", "!/time/$id") | Out-Null
}

# Hyperlinks.Add() re-stamps the cell style; put the original "Hyperlink" cell
# style back on every link cell so formatting matches the rest of the column.
$ws.Range("C2:C21").Style = "Hyperlink"
